# Fixed Tests for SamplePatholoy, SampleType, Se, StageOfDisease, and Study
#
# The "CasesTab" query (cell B2 on the startup sheet) is updated to drop the
# trailing `Cohort` column from its Cypher RETURN clause.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newCasesTabQuery = @"
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
MATCH (samp:sample)-->(c) 
 WHERE samp.specific_sample_pathology IN ["T Cell Lymphoma"]  
OPTIONAL MATCH (co:cohort)<-[*]-(c)
  WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS ``Case ID`` ,
        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,
        coalesce(s.clinical_study_type, '') AS  ``Study Type``,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,
        coalesce(demo.weight, '') AS ``Weight (kg)``,
        coalesce(diag.best_response, '') AS ``Response to Treatment``
"@

# Trim the trailing newline the here-string introduces before the closing "@
$newCasesTabQuery = $newCasesTabQuery.TrimEnd("`r", "`n")

$ws.Range("B2").Value = $newCasesTabQuery

# Excel's own re-save moved the active selection from B4 to B2 as part of this
# edit - mirror that so the saved sheetView reflects it too.
$ws.Range("B2").Select()

Write-Output "Updated startup!B2 CasesTab query (removed Cohort column)."
